# Birdoptera pace calculator - "Changed miles to Miles/Km"
# Shortens Aid-station names and relabels the miles/km headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row updates -------------------------------------------------
$ws.Range("I1").Value = "Min per Miles/Km ="
$ws.Range("B2").Value = "Miles/ Km"
$ws.Range("C2").Value = "Total Miles/ Km"

# --- Aid station names (column A, rows 3-21) -----------------------------
$ws.Range("A3").Value = "START "
$ws.Range("A4").Value = "Aid #1"
$ws.Range("A5").Value = "Aid #2"
$ws.Range("A6").Value = "Aid #3"
$ws.Range("A7").Value = "Aid #4"
$ws.Range("A8").Value = "Aid #5"
$ws.Range("A9").Value = "Aid #6"
$ws.Range("A10").Value = "Aid #7"
$ws.Range("A11").Value = "Aid #8"
$ws.Range("A12").Value = "Aid #9"
$ws.Range("A13").Value = "Aid #10"
$ws.Range("A14").Value = "Aid #11"
$ws.Range("A15").Value = "Aid #12"
$ws.Range("A16").Value = "Aid #13"
$ws.Range("A17").Value = "Aid #14"
$ws.Range("A18").Value = "Aid #15"
$ws.Range("A19").Value = "Aid #16"
$ws.Range("A20").Value = "Aid #17"
$ws.Range("A21").Value = "Finish!"

# --- Selection cosmetics (matches the saved file's cursor position) -----
$ws.Range("A1:A2").Select()
